$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for team record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, bordered, centered) from an existing header cell
# onto the new header cells so they match the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team record (constant for every player row): 95 wins, 67 losses, 0 ties
$ws.Range("AD2:AD54").Value = 95
$ws.Range("AE2:AE54").Value = 67
$ws.Range("AF2:AF54").Value = 0
